$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$teams = @(
    "ATALANTA",
    "BOLOGNA",
    "CREMONESE",
    "EMPOLI",
    "FIORENTINA",
    "INTER",
    "JUVENTUS",
    "LAZIO",
    "LECCE",
    "MILAN",
    "MONZA",
    "NAPOLI",
    "ROMA",
    "SALERNITANA",
    "SAMPDORIA",
    "SASSUOLO",
    "SPEZIA",
    "TORINO",
    "UDINESE",
    "VERONA"
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $teams[$i]
}
